$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J2").Value = 3.1
$ws.Range("M2").Value = 1.05
$ws.Range("O2").Value = 1.37
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.67
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("AB2").Value = 41
$ws.Range("AE2").Value = 19
$ws.Range("AH2").Value = 9
$ws.Range("AP2").Value = 26
$ws.Range("AT2").Value = 2.5
$ws.Range("BB2").Value = 251
$ws.Range("G3").Value = 2.38
$ws.Range("I3").Value = 3.4
$ws.Range("M3").Value = 1.07
$ws.Range("O3").Value = 1.37
$ws.Range("Q3").Value = 2.3
$ws.Range("R3").Value = 1.62
$ws.Range("AH3").Value = 9
$ws.Range("BA3").Value = 101
$ws.Range("G4").Value = 1.9
$ws.Range("H4").Value = 3.3
$ws.Range("I4").Value = 4.33
$ws.Range("J4").Value = 2.63
$ws.Range("L4").Value = 4.75
$ws.Range("M4").Value = 1.05
$ws.Range("O4").Value = 1.37
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.3
$ws.Range("R4").Value = 1.62
$ws.Range("U4").Value = 2.05
$ws.Range("V4").Value = 1.7
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 8
$ws.Range("Z4").Value = 15
$ws.Range("AG4").Value = 451
$ws.Range("AH4").Value = 10
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 15
$ws.Range("AM4").Value = 41
$ws.Range("AN4").Value = 3.75
$ws.Range("AO4").Value = 11
$ws.Range("AQ4").Value = 41
$ws.Range("AW4").Value = 6
$ws.Range("AX4").Value = 26
$ws.Range("AY4").Value = 34
$ws.Range("AZ4").Value = 81
$ws.Range("BA4").Value = 126
$ws.Range("BB4").Value = 301
$ws.Range("BD4").Value = 126
$ws.Range("V5").Value = 1.63
$ws.Range("Q8").Value = 1.95
$ws.Range("R8").Value = 1.9
$ws.Range("U8").Value = 1.77
$ws.Range("V8").Value = 1.87
$ws.Range("M23").Value = 1.08
$ws.Range("O23").Value = 1.44
$ws.Range("P23").Value = 2.63
$ws.Range("R23").Value = 1.54
$ws.Range("V23").Value = 1.73
$ws.Range("M24").Value = 1.08
$ws.Range("O24").Value = 1.4
$ws.Range("R24").Value = 1.57
$ws.Range("U24").Value = 1.91
$ws.Range("V24").Value = 1.8
$ws.Range("M30").Value = 1.07
$ws.Range("O30").Value = 1.3
$ws.Range("R30").Value = 1.8
$ws.Range("S30").Value = 1.4
$ws.Range("T30").Value = 2.75
$ws.Range("X30").Value = 10
$ws.Range("AT30").Value = 2.75
$ws.Range("AX30").Value = 19
$ws.Range("M38").Value = 1.05
$ws.Range("O38").Value = 1.41
$ws.Range("P38").Value = 2.62
$ws.Range("R38").Value = 1.5
$ws.Range("M39").Value = 1.03
$ws.Range("O39").Value = 1.25
